# Apply weekly crypto price/volume refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.027.08'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -3.13%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.187.70'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.15%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.28'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -2.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '169.47'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -6.34%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.612'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -5.92%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.189.74'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -3.03%  '

# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -3.20%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.74'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.19%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.388'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -4.41%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.741.71'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.16%  '

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.76%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '64.156.47'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.98%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.52'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -3.48%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.26%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.177.52'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -3.81%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '416.69'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -3.96%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.39'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.08%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.88'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.91%  '

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.97%  '

# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.08%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.79'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.11%  '

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +3.27%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.496'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -3.11%  '

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -3.08%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.57%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.10%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.86'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.85%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.03'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.45%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.999'

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.04'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.84%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.43'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -2.98%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '156.27'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.42%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.39'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.89%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.748.03'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.18%  '

# Row 39
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.85%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.21'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -5.10%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.21'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.04%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.720'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -6.68%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.92'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -3.09%  '

# Row 44
$ws.Range("B44").Value = 'RenderToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.74'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.71%  '

# Row 45
$ws.Range("B45").Value = 'Hedera'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0631'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.39%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '22.19'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -4.72%  '

# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0264'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.63%  '

# Row 48
$ws.Range("B48").Value = 'Bittensor'
$ws.Range("C48").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '298.72'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.67%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.11'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -8.88%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0998'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -5.28%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.02%  '
